$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 218.6
$ws.Range("I9").Value = 198
$ws.Range("J9").Value = 249.5
$ws.Range("K9").Value = 198
$ws.Range("L9").Value = 249.5
$ws.Range("M9").Value = -29
$ws.Range("N9").Value = -587.5
$ws.Range("H17").Value = 2966.4424
$ws.Range("J17").Value = 3144.8958
$ws.Range("L17").Value = 9434.687399999999
$ws.Range("N17").Value = -9770.687399999999
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("H54").Value = 502000
$ws.Range("I54").Value = 502000
$ws.Range("K54").Value = 502000
$ws.Range("M54").Value = -501514
$ws.Range("H106").Value = 7148.0586
$ws.Range("I106").Value = 2364.6365
$ws.Range("K106").Value = 2364.6365
$ws.Range("M106").Value = -1733.6365
$ws.Range("H137").Value = 3103.9707
$ws.Range("I137").Value = 2344.125
$ws.Range("J137").Value = 3779.389
$ws.Range("K137").Value = 7032.375
$ws.Range("L137").Value = 11338.167
$ws.Range("M137").Value = -4482.375
$ws.Range("N137").Value = -16438.167
$ws.Range("H138").Value = 4130.1924
$ws.Range("J138").Value = 5499.6
$ws.Range("L138").Value = 16498.8
$ws.Range("N138").Value = -26778.8
$ws.Range("H141").Value = 959.2857
$ws.Range("I141").Value = 959.2857
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 2877.8571
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2302.1429
$ws.Range("N26").Value = ""
$ws.Range("N141").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 8767.083000000001
$ws.Range("I2").Value = 473.18182
$ws.Range("K2").Value = 473.18182
$ws.Range("M2").Value = -360.18182
$ws.Range("H37").Value = 10259.5
$ws.Range("I37").Value = 3666.6667
$ws.Range("K37").Value = 3666.6667
$ws.Range("M37").Value = -3393.6667
$ws.Range("H43").Value = 31981.111
$ws.Range("J43").Value = 31355.428
$ws.Range("L43").Value = 31355.428
$ws.Range("N43").Value = -31981.428
$ws.Range("H45").Value = 90912370
$ws.Range("J45").Value = 6596
$ws.Range("L45").Value = 6596
$ws.Range("N45").Value = -7350
$ws.Range("H116").Value = 8767.083000000001
$ws.Range("I116").Value = 473.18182
$ws.Range("K116").Value = 473.18182
$ws.Range("M116").Value = 1820.81818

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 8767.083000000001
$ws.Range("I3").Value = 473.18182
$ws.Range("K3").Value = 473.18182
$ws.Range("M3").Value = -359.18182
$ws.Range("H86").Value = 3588.5417
$ws.Range("I86").Value = 2911.9
$ws.Range("J86").Value = 6971.75
$ws.Range("K86").Value = 2911.9
$ws.Range("L86").Value = 6971.75
$ws.Range("M86").Value = -1788.9
$ws.Range("N86").Value = -9217.75
$ws.Range("H89").Value = 3588.5417
$ws.Range("I89").Value = 2911.9
$ws.Range("J89").Value = 6971.75
$ws.Range("K89").Value = 14559.5
$ws.Range("L89").Value = 34858.75
$ws.Range("M89").Value = -8943.5
$ws.Range("N89").Value = -46090.75
$ws.Range("H107").Value = 2344.0908
$ws.Range("I107").Value = 2078.5
$ws.Range("K107").Value = 2078.5
$ws.Range("M107").Value = -158.5
$ws.Range("H134").Value = 1544.6786
$ws.Range("I134").Value = 1045.8148
$ws.Range("K134").Value = 3137.4444
$ws.Range("M134").Value = -602.4444000000003

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32607.805
$ws.Range("I31").Value = 3386.8518
$ws.Range("J31").Value = 120270.664
$ws.Range("K31").Value = 3386.8518
$ws.Range("L31").Value = 120270.664
$ws.Range("M31").Value = -3091.8518
$ws.Range("N31").Value = -120860.664
$ws.Range("H34").Value = 32607.805
$ws.Range("I34").Value = 3386.8518
$ws.Range("J34").Value = 120270.664
$ws.Range("K34").Value = 3386.8518
$ws.Range("L34").Value = 120270.664
$ws.Range("M34").Value = -3184.8518
$ws.Range("N34").Value = -120674.664
$ws.Range("H58").Value = 3285
$ws.Range("I58").Value = 1674.0769
$ws.Range("K58").Value = 1674.0769
$ws.Range("M58").Value = -1471.0769
$ws.Range("H105").Value = 8499.777
$ws.Range("I105").Value = 7185.875
$ws.Range("K105").Value = 7185.875
$ws.Range("M105").Value = -5438.875
$ws.Range("H107").Value = 1267.3158
$ws.Range("I107").Value = 1524.6
$ws.Range("K107").Value = 1524.6
$ws.Range("M107").Value = 395.4000000000001
$ws.Range("H134").Value = 2673.4546
$ws.Range("I134").Value = 1517.7333
$ws.Range("J134").Value = 5150
$ws.Range("K134").Value = 4553.199900000001
$ws.Range("L134").Value = 15450
$ws.Range("M134").Value = -2018.199900000001
$ws.Range("N134").Value = -20520
$ws.Range("H136").Value = 3285
$ws.Range("I136").Value = 1674.0769
$ws.Range("K136").Value = 5022.2307
$ws.Range("M136").Value = -2472.2307

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 46512.92
$ws.Range("I2").Value = 97
$ws.Range("J2").Value = 100664.836
$ws.Range("K2").Value = 582
$ws.Range("L2").Value = 603989.0159999999
$ws.Range("M2").Value = -469
$ws.Range("N2").Value = -604215.0159999999
$ws.Range("H23").Value = 367.57144
$ws.Range("J23").Value = 395.5
$ws.Range("L23").Value = 1186.5
$ws.Range("N23").Value = -1656.5
$ws.Range("H42").Value = 11334.667
$ws.Range("J42").Value = 14502
$ws.Range("L42").Value = 43506
$ws.Range("N42").Value = -44574
$ws.Range("H75").Value = 76928340
$ws.Range("I75").Value = 250000480
$ws.Range("J75").Value = 7388
$ws.Range("K75").Value = 750001440
$ws.Range("L75").Value = 22164
$ws.Range("M75").Value = -750000442
$ws.Range("N75").Value = -24160
$ws.Range("H78").Value = 76928340
$ws.Range("I78").Value = 250000480
$ws.Range("J78").Value = 7388
$ws.Range("K78").Value = 2250004320
$ws.Range("L78").Value = 66492
$ws.Range("M78").Value = -2249999328
$ws.Range("N78").Value = -76476
$ws.Range("H80").Value = 10167.667
$ws.Range("I80").Value = 8001
$ws.Range("K80").Value = 24003
$ws.Range("M80").Value = -23067
$ws.Range("H83").Value = 10167.667
$ws.Range("I83").Value = 8001
$ws.Range("K83").Value = 72009
$ws.Range("M83").Value = -67329
$ws.Range("H86").Value = 1913
$ws.Range("J86").Value = 3799.5
$ws.Range("L86").Value = 11398.5
$ws.Range("N86").Value = -13770.5
$ws.Range("H89").Value = 1913
$ws.Range("J89").Value = 3799.5
$ws.Range("L89").Value = 34195.5
$ws.Range("N89").Value = -46051.5
$ws.Range("H131").Value = 35496052
$ws.Range("I131").Value = 41670936
$ws.Range("J131").Value = 32408608
$ws.Range("K131").Value = 125012808
$ws.Range("L131").Value = 97225824
$ws.Range("M131").Value = -125007768
$ws.Range("N131").Value = -97235904
$ws.Range("H139").Value = 5449.25
$ws.Range("I139").Value = 1608.4286
$ws.Range("J139").Value = 10826.4
$ws.Range("K139").Value = 4825.2858
$ws.Range("L139").Value = 32479.2
$ws.Range("M139").Value = 314.7142000000003
$ws.Range("N139").Value = -42759.2
$ws.Range("H140").Value = 3700.1875
$ws.Range("I140").Value = 2862.6155
$ws.Range("K140").Value = 8587.8465
$ws.Range("M140").Value = -3407.8465

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 36021
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("H95").Value = 15000
$ws.Range("J95").Value = 15000
$ws.Range("L95").Value = 15000
$ws.Range("H113").Value = 4235.2
$ws.Range("I113").Value = 3652.2104
$ws.Range("J113").Value = 6081.3335
$ws.Range("K113").Value = 3652.2104
$ws.Range("L113").Value = 6081.3335
$ws.Range("M113").Value = -1482.2104
$ws.Range("N113").Value = -10421.3335
$ws.Range("H132").Value = 23465.22
$ws.Range("I132").Value = 26882.334
$ws.Range("K132").Value = 80647.00199999999
$ws.Range("M132").Value = -78117.00199999999
$ws.Range("M53").Value = ""
$ws.Range("N95").Value = -20492

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5537.4116
$ws.Range("J22").Value = 7649.1816
$ws.Range("L22").Value = 7649.1816
$ws.Range("N22").Value = -8239.1816
$ws.Range("H27").Value = 5537.4116
$ws.Range("J27").Value = 7649.1816
$ws.Range("L27").Value = 7649.1816
$ws.Range("N27").Value = -7863.1816
$ws.Range("H46").Value = 2596.7144
$ws.Range("J46").Value = 2727.2307
$ws.Range("L46").Value = 2727.2307
$ws.Range("N46").Value = -3103.2307
$ws.Range("H55").Value = 2382312
$ws.Range("I55").Value = 4545792
$ws.Range("K55").Value = 4545792
$ws.Range("M55").Value = -4545619
$ws.Range("H132").Value = 3356
$ws.Range("I132").Value = 1647.2222
$ws.Range("K132").Value = 4941.6666
$ws.Range("M132").Value = -2411.6666
$ws.Range("H136").Value = 12763.167
$ws.Range("I136").Value = 2787.5
$ws.Range("K136").Value = 8362.5
$ws.Range("M136").Value = -5812.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("H132").Value = 5133.3184
$ws.Range("I132").Value = 4296.4
$ws.Range("K132").Value = 12889.2
$ws.Range("M132").Value = -10359.2
$ws.Range("H136").Value = 5578.15
$ws.Range("I136").Value = 2840.9285
$ws.Range("J136").Value = 11965
$ws.Range("K136").Value = 8522.7855
$ws.Range("L136").Value = 35895
$ws.Range("M136").Value = -5972.7855
$ws.Range("N136").Value = -40995
$ws.Range("N110").Value = ""
